$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column C (header "Förändrad") holds a date serial number (45179 -> 2023-09-10).
# Update every data row (2..135) from 45179 to 45180 (2023-09-11), keeping the
# existing number formatting/style of the cells untouched.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 3).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45179) {
        $cell.Value2 = 45180
    }
}
